$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: the stray "_GoBack" bookmark that used to sit between "Multiple
# pairs of " and "either the black..." is gone (the user's cursor moved on),
# so those two text fragments collapse back into a single contiguous run.
# ---------------------------------------------------------------------------

$mergedSentence = "Multiple pairs of either the black or brown could be drawn before a white pair was drawn.  Since there are only 4 white socks total and 2 of them need to be selected for a pair, all of the colored socks plus 2 of the white could potentially be selected prior to getting a matching white pair.  "

$findRange = $d.Content
$found = $findRange.Find.Execute($mergedSentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$spanStart = $findRange.Start
$spanEnd = $findRange.End

$toDelete = $d.Range($spanStart, $spanEnd)
$toDelete.Delete()

$reinsert = $d.Range($spanStart, $spanStart)
$reinsert.InsertAfter($mergedSentence)

# ---------------------------------------------------------------------------
# Part 2: add the new answer text for question 5a (Socks in the Dark ->
# "Choose a solution and develop a plan to implement it:") right after the
# existing "a. " marker, followed by a second paragraph ("To select a
# matching pair for each color...") whose _GoBack bookmark marks the spot
# right before "its match.".
# ---------------------------------------------------------------------------

$header = "Choose a solution and develop a plan to implement it:"

# Locate the first occurrence (Cat/Parrot problem) so we can search past it
# and land on the second occurrence that belongs to the Socks problem.
$headerRange = $d.Content
$null = $headerRange.Find.Execute($header, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$secondHeaderRange = $d.Content
$secondHeaderRange.Start = $headerRange.End
$null = $secondHeaderRange.Find.Execute($header, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Now find the "a. " marker that immediately follows that second header.
$aMarkerRange = $d.Content
$aMarkerRange.Start = $secondHeaderRange.End
$null = $aMarkerRange.Find.Execute("a. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPos = $aMarkerRange.End

$part1Text = "To select at least one matching pair: all of the socks would be combined in a drawer.  A person would draw one sock.  To be able to guarantee a match, a total of 3 more socks would need to be drawn.  It is possible that fewer than 3 would need to be drawn, but by drawing 3 more it is guaranteed that one of those 3 additional socks would be a match for the first one."
$lineBreak = [char]11
$part2aText = "To select a matching pair for each color: all of the socks would be combined in a drawer.  A person would draw one sock and then 3 more to guarantee "
$part2bText = "its match.  The same would be done for each additional color.  Because there are only 4 white socks, it is possible that all 16 colored socks could be drawn before the 4 white socks are left.  Therefore in order to guarantee a match of each color, a minimum of 18 socks would need to be drawn.  Again, it is possible that fewer than 18 would solve this problem, but it is not a guarantee.  "

$fullBlock = $part1Text + $lineBreak + $lineBreak + $part2aText + $part2bText

$insertion = $d.Range($insertPos, $insertPos)
$insertion.InsertAfter($fullBlock)

# Re-seat the "_GoBack" bookmark in its new, final resting place: right
# before "its match." (i.e., right after "...to guarantee ").
$bookmarkPos = $insertPos + $part1Text.Length + 2 + $part2aText.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
